$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.059611
$ws.Range("H2").Value = 0.178833
$ws.Range("I2").Value = 0.003943014985542741
$ws.Range("J2").Value = 0.003943014985542741
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.949891
$ws.Range("N2").Value = 2.849673
$ws.Range("O2").Value = 0.1664393778377885
$ws.Range("P2").Value = 0.1664393778377885
$ws.Range("Q2").Value = 0.056623952401
$ws.Range("R2").Value = 0.509615571609
$ws.Range("S2").Value = 0.0006562729609988104
$ws.Range("T2").Value = 0.0006562729609988104
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.059611
$ws.Range("H3").Value = 0.178833
$ws.Range("I3").Value = 0.003943014985542741
$ws.Range("J3").Value = 0.003943014985542741
$ws.Range("O3").Value = 0.4967272219242518
$ws.Range("P3").Value = 0.4967272219242519
$ws.Range("Q3").Value = 0.1689904092163333
$ws.Range("R3").Value = 1.520913682947
$ws.Range("S3").Value = 0.00195860287977434
$ws.Range("T3").Value = 0.00195860287977434
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.059611
$ws.Range("H4").Value = 0.178833
$ws.Range("I4").Value = 0.003943014985542741
$ws.Range("J4").Value = 0.003943014985542741
$ws.Range("M4").Value = 1.722217666666667
$ws.Range("N4").Value = 5.166653
$ws.Range("O4").Value = 0.3017660309880268
$ws.Range("P4").Value = 0.3017660309880269
$ws.Range("Q4").Value = 0.1026631173276667
$ws.Range("R4").Value = 0.923968055949
$ws.Range("S4").Value = 0.001189867982313545
$ws.Range("T4").Value = 0.001189867982313545
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.059611
$ws.Range("H5").Value = 0.178833
$ws.Range("I5").Value = 0.003943014985542741
$ws.Range("J5").Value = 0.003943014985542741
$ws.Range("M5").Value = 0.200134
$ws.Range("N5").Value = 0.600402
$ws.Range("O5").Value = 0.03506736924993285
$ws.Range("P5").Value = 0.03506736924993285
$ws.Range("Q5").Value = 0.011930187874
$ws.Range("R5").Value = 0.107371690866
$ws.Range("S5").Value = 0.0001382711624560459
$ws.Range("T5").Value = 0.0001382711624560459
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.75172666666667
$ws.Range("H6").Value = 44.25518
$ws.Range("I6").Value = 0.9757641930062764
$ws.Range("J6").Value = 0.9757641930062765
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.949891
$ws.Range("N6").Value = 2.849673
$ws.Range("O6").Value = 0.1664393778377885
$ws.Range("P6").Value = 0.1664393778377885
$ws.Range("Q6").Value = 14.01253239512667
$ws.Range("R6").Value = 126.11279155614
$ws.Range("S6").Value = 0.1624055852003564
$ws.Range("T6").Value = 0.1624055852003564
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 14.75172666666667
$ws.Range("H7").Value = 44.25518
$ws.Range("I7").Value = 0.9757641930062764
$ws.Range("J7").Value = 0.9757641930062765
$ws.Range("O7").Value = 0.4967272219242518
$ws.Range("P7").Value = 0.4967272219242519
$ws.Range("Q7").Value = 41.81946832040223
$ws.Range("R7").Value = 376.3752148836201
$ws.Range("S7").Value = 0.4846886368451672
$ws.Range("T7").Value = 0.4846886368451673
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.75172666666667
$ws.Range("H8").Value = 44.25518
$ws.Range("I8").Value = 0.9757641930062764
$ws.Range("J8").Value = 0.9757641930062765
$ws.Range("M8").Value = 1.722217666666667
$ws.Range("N8").Value = 5.166653
$ws.Range("O8").Value = 0.3017660309880268
$ws.Range("P8").Value = 0.3017660309880269
$ws.Range("Q8").Value = 25.40568427917111
$ws.Range("R8").Value = 228.65115851254
$ws.Range("S8").Value = 0.294452487703739
$ws.Range("T8").Value = 0.2944524877037391
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.75172666666667
$ws.Range("H9").Value = 44.25518
$ws.Range("I9").Value = 0.9757641930062764
$ws.Range("J9").Value = 0.9757641930062765
$ws.Range("M9").Value = 0.200134
$ws.Range("N9").Value = 0.600402
$ws.Range("O9").Value = 0.03506736924993285
$ws.Range("P9").Value = 0.03506736924993285
$ws.Range("Q9").Value = 2.952322064706667
$ws.Range("R9").Value = 26.57089858236
$ws.Range("S9").Value = 0.03421748325701383
$ws.Range("T9").Value = 0.03421748325701384
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.306789
$ws.Range("H10").Value = 0.9203669999999999
$ws.Range("I10").Value = 0.0202927920081809
$ws.Range("J10").Value = 0.02029279200818091
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.949891
$ws.Range("N10").Value = 2.849673
$ws.Range("O10").Value = 0.1664393778377885
$ws.Range("P10").Value = 0.1664393778377885
$ws.Range("Q10").Value = 0.291416109999
$ws.Range("R10").Value = 2.622744989991
$ws.Range("S10").Value = 0.003377519676433276
$ws.Range("T10").Value = 0.003377519676433277
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.306789
$ws.Range("H11").Value = 0.9203669999999999
$ws.Range("I11").Value = 0.0202927920081809
$ws.Range("J11").Value = 0.02029279200818091
$ws.Range("O11").Value = 0.4967272219242518
$ws.Range("P11").Value = 0.4967272219242519
$ws.Range("Q11").Value = 0.869711943317
$ws.Range("R11").Value = 7.827407489853
$ws.Range("S11").Value = 0.01007998219931036
$ws.Range("T11").Value = 0.01007998219931036
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.306789
$ws.Range("H12").Value = 0.9203669999999999
$ws.Range("I12").Value = 0.0202927920081809
$ws.Range("J12").Value = 0.02029279200818091
$ws.Range("M12").Value = 1.722217666666667
$ws.Range("N12").Value = 5.166653
$ws.Range("O12").Value = 0.3017660309880268
$ws.Range("P12").Value = 0.3017660309880269
$ws.Range("Q12").Value = 0.528357435739
$ws.Range("R12").Value = 4.755216921651
$ws.Range("S12").Value = 0.006123675301974303
$ws.Range("T12").Value = 0.006123675301974304
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.306789
$ws.Range("H13").Value = 0.9203669999999999
$ws.Range("I13").Value = 0.0202927920081809
$ws.Range("J13").Value = 0.02029279200818091
$ws.Range("M13").Value = 0.200134
$ws.Range("N13").Value = 0.600402
$ws.Range("O13").Value = 0.03506736924993285
$ws.Range("P13").Value = 0.03506736924993285
$ws.Range("Q13").Value = 0.061398909726
$ws.Range("R13").Value = 0.5525901875339999
$ws.Range("S13").Value = 0.000711614830462966
$ws.Range("T13").Value = 0.0007116148304629661
